# Generate Report for Handoff
# Updates the localization-status report to reflect that the content is
# now "Ready for handoff" (previously "In Translation"), refreshes the
# associated timestamps, and widens the status columns so the new,
# longer status text fits (mirrors an auto-fit of the Status columns).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-08-30 19:13:34"
$ws.Columns.Item(5).ColumnWidth = 16.33
$ws.Columns.Item(6).ColumnWidth = 16.33

# ---- zh-cn sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-30 19:13:29"
$ws.Columns.Item(3).ColumnWidth = 16.33

# ---- de-de sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-30 19:13:34"
$ws.Columns.Item(3).ColumnWidth = 16.33
